$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.144.11"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "1.834.50"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'240.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "'0.6648"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.52%  "

$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.2956"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.77%  "

$ws.Range("D9").Value = "'0.07358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.34%  "

$ws.Range("D10").Value = "'22.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.54%  "

$ws.Range("D11").Value = "'0.07679"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "

$ws.Range("D12").Value = "1.835.89"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").Value = "'5.025"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.65%  "

$ws.Range("D14").Value = "'0.6751"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("D15").Value = "'86.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.24%  "

$ws.Range("D16").Value = "'6.169"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "29.140.69"
$ws.Range("E17").Value = "  -1.14%  "

$ws.Range("D18").Value = "'0.000008244"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").Value = "'229.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.60%  "

$ws.Range("D20").Value = "'12.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.74%  "

$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'7.300"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'0.9994"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'160.70"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.1422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.92%  "

$ws.Range("D26").Value = "'8.672"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "

$ws.Range("D27").Value = "'18.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'1.505"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("D29").Value = "'4.237"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("D30").Value = "'4.096"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("D31").Value = "'1.202"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").Value = "'0.05327"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "

$ws.Range("E33").Value = "  -1.09%  "

$ws.Range("D34").Value = "'0.7457"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("D36").Value = "'2.677"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("D37").Value = "1.318.43"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "'0.01806"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.62%  "

$ws.Range("D39").Value = "'2.713"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").Value = "'0.9254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.75%  "

$ws.Range("D41").Value = "'6.010"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.96%  "

$ws.Range("D42").Value = "'0.9981"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").Value = "'103.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "

$ws.Range("D44").Value = "1.983.02"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").Value = "'0.5168"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("E46").Value = "  -3.08%  "

$ws.Range("D47").Value = "'1.761"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Value = "'63.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("D49").Value = "'9.279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.53%  "

$ws.Range("D50").Value = "'0.07469"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.47%  "

$ws.Range("D51").Value = "'0.05925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
